$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the two new columns
$ws.Range("AB1").Value = "mean_transitivity"
$ws.Range("AC1").Value = "sd_transitivity"

# Fill in the computed mean/sd transitivity values per word row
$ws.Range("AB2").Value = 0.05689396431337302
$ws.Range("AC2").Value = 0.11376179957922743
$ws.Range("AB3").Value = 0.36794537363817703
$ws.Range("AC3").Value = 0.3896714219499534
$ws.Range("AB4").Value = 0.48032128514056227
$ws.Range("AC4").Value = 0.42301781163298774
$ws.Range("AB5").Value = 0.44828869047619047
$ws.Range("AC5").Value = 0.461540222225728
$ws.Range("AB6").Value = 0.34202898550724636
$ws.Range("AC6").Value = 0.42716775775464627
$ws.Range("AB7").Value = 0.38606965174129354
$ws.Range("AC7").Value = 0.45342502243191796
$ws.Range("AB8").Value = 0.23532743851892787
$ws.Range("AC8").Value = 0.3683923374898022
$ws.Range("AB9").Value = 0.24996488270824552
$ws.Range("AC9").Value = 0.35955910393658086
$ws.Range("AB10").Value = 0.35035014005602244
$ws.Range("AC10").Value = 0.4175583442395314
$ws.Range("AB11").Value = 0.3036231884057971
$ws.Range("AC11").Value = 0.3971958031316237
$ws.Range("AB12").Value = 0.3371087375760273
$ws.Range("AC12").Value = 0.40633481417078166
$ws.Range("AB13").Value = 0.40534722222222225
$ws.Range("AC13").Value = 0.4412769014557637
$ws.Range("AB14").Value = 0.37331349206349207
$ws.Range("AC14").Value = 0.3910193603746583
$ws.Range("AB15").Value = 0.21734052111410598
$ws.Range("AC15").Value = 0.32215198666431
$ws.Range("AB16").Value = 0.34193986928104575
$ws.Range("AC16").Value = 0.3639908132010843
$ws.Range("AB17").Value = 0.24537318473488684
$ws.Range("AC17").Value = 0.3310940322401557
$ws.Range("AB18").Value = 0.4094062316284538
$ws.Range("AC18").Value = 0.42330559614300817
$ws.Range("AB19").Value = 0.27559523809523806
$ws.Range("AC19").Value = 0.3933194950170748
$ws.Range("AB20").Value = 0.2203220015004912
$ws.Range("AC20").Value = 0.3017504751563466
$ws.Range("AB21").Value = 0.41690821256038646
$ws.Range("AC21").Value = 0.4494958391631963
$ws.Range("AB22").Value = 0.35654761904761906
$ws.Range("AC22").Value = 0.4316713884312784
$ws.Range("AB23").Value = 0.3044248749728202
$ws.Range("AC23").Value = 0.41298016256242365
$ws.Range("AB24").Value = 0.37918706423379317
$ws.Range("AC24").Value = 0.3806844499983209
$ws.Range("AB25").Value = 0.36614583333333334
$ws.Range("AC25").Value = 0.40751874281711076
$ws.Range("AB26").Value = 0.43056668150126093
$ws.Range("AC26").Value = 0.41648583997433053
$ws.Range("AB27").Value = 0.24793522682475008
$ws.Range("AC27").Value = 0.3440469473938018
$ws.Range("AB28").Value = 0.23640444522797463
$ws.Range("AC28").Value = 0.30460522732632617
$ws.Range("AB29").Value = 0.28593576965669987
$ws.Range("AC29").Value = 0.3677085973313867
$ws.Range("AB30").Value = 0.35534373882997733
$ws.Range("AC30").Value = 0.391412096267067
$ws.Range("AB31").Value = 0.4497777777777778
$ws.Range("AC31").Value = 0.43929965192813913
$ws.Range("AB32").Value = 0.3973356009070295
$ws.Range("AC32").Value = 0.4081702210747002
$ws.Range("AB33").Value = 0.3553375196232339
$ws.Range("AC33").Value = 0.4022070537253573
$ws.Range("AB34").Value = 0.37448634444047285
$ws.Range("AC34").Value = 0.4096735762357971
